$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4844.654
$ws.Range("I70").Value = 1959.25
$ws.Range("J70").Value = 6127.0557
$ws.Range("K70").Value = 5877.75
$ws.Range("L70").Value = 18381.1671
$ws.Range("M70").Value = -5607.75
$ws.Range("N70").Value = -18921.1671
$ws.Range("H73").Value = 4844.654
$ws.Range("I73").Value = 1959.25
$ws.Range("J73").Value = 6127.0557
$ws.Range("K73").Value = 5877.75
$ws.Range("L73").Value = 18381.1671
$ws.Range("M73").Value = -4941.75
$ws.Range("N73").Value = -20253.1671
$ws.Range("H76").Value = 4883.8857
$ws.Range("J76").Value = 5669.8
$ws.Range("L76").Value = 5669.8
$ws.Range("N76").Value = -6299.8
$ws.Range("H79").Value = 4883.8857
$ws.Range("J79").Value = 5669.8
$ws.Range("L79").Value = 5669.8
$ws.Range("N79").Value = -7853.8
$ws.Range("H106").Value = 7665.8335
$ws.Range("J106").Value = 8699
$ws.Range("L106").Value = 8699
$ws.Range("N106").Value = -9961
$ws.Range("H137").Value = 8182.14
$ws.Range("I137").Value = 11371.226
$ws.Range("J137").Value = 2978.8948
$ws.Range("K137").Value = 34113.678
$ws.Range("L137").Value = 8936.6844
$ws.Range("M137").Value = -31563.678
$ws.Range("N137").Value = -14036.6844
$ws.Range("H138").Value = 3431.9412
$ws.Range("I138").Value = 2354.111
$ws.Range("J138").Value = 4644.5
$ws.Range("K138").Value = 7062.333
$ws.Range("L138").Value = 13933.5
$ws.Range("M138").Value = -1922.333
$ws.Range("N138").Value = -24213.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1672.1666
$ws.Range("I2").Value = 1368.4584
$ws.Range("K2").Value = 1368.4584
$ws.Range("M2").Value = -1255.4584
$ws.Range("H45").Value = 4375.15
$ws.Range("J45").Value = 5265.2144
$ws.Range("L45").Value = 5265.2144
$ws.Range("N45").Value = -6019.2144
$ws.Range("H61").Value = 2970.8
$ws.Range("I61").Value = 2039.6522
$ws.Range("K61").Value = 2039.6522
$ws.Range("M61").Value = -1827.6522
$ws.Range("H97").Value = 2896.4443
$ws.Range("I97").Value = 2472.125
$ws.Range("J97").Value = 3513.6365
$ws.Range("K97").Value = 2472.125
$ws.Range("L97").Value = 3513.6365
$ws.Range("M97").Value = -1976.125
$ws.Range("N97").Value = -4505.636500000001
$ws.Range("H102").Value = 2994.6667
$ws.Range("I102").Value = 3039.6365
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 3039.6365
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -1417.6365
$ws.Range("N102").Value = -5744
$ws.Range("H116").Value = 1672.1666
$ws.Range("I116").Value = 1368.4584
$ws.Range("K116").Value = 1368.4584
$ws.Range("M116").Value = 925.5416
$ws.Range("H136").Value = 2970.8
$ws.Range("I136").Value = 2039.6522
$ws.Range("K136").Value = 6118.9566
$ws.Range("M136").Value = -3568.9566

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1672.1666
$ws.Range("I3").Value = 1368.4584
$ws.Range("K3").Value = 1368.4584
$ws.Range("M3").Value = -1254.4584
$ws.Range("H86").Value = 1027.027
$ws.Range("I86").Value = 1026.7
$ws.Range("K86").Value = 1026.7
$ws.Range("M86").Value = 96.29999999999995
$ws.Range("H89").Value = 1027.027
$ws.Range("I89").Value = 1026.7
$ws.Range("K89").Value = 5133.5
$ws.Range("M89").Value = 482.5
$ws.Range("H94").Value = 998.4074000000001
$ws.Range("I94").Value = 1087.0555
$ws.Range("J94").Value = 821.1111
$ws.Range("K94").Value = 1087.0555
$ws.Range("L94").Value = 821.1111
$ws.Range("M94").Value = -636.0554999999999
$ws.Range("N94").Value = -1723.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 54384.418
$ws.Range("I86").Value = 68845.11
$ws.Range("K86").Value = 68845.11
$ws.Range("M86").Value = -67722.11
$ws.Range("H89").Value = 54384.418
$ws.Range("I89").Value = 68845.11
$ws.Range("K89").Value = 344225.55
$ws.Range("M89").Value = -338609.55
$ws.Range("H132").Value = 27547.521
$ws.Range("I132").Value = 29694.953
$ws.Range("K132").Value = 89084.859
$ws.Range("M132").Value = -86554.859
$ws.Range("H135").Value = 119123.75
$ws.Range("J135").Value = 119123.75
$ws.Range("L135").Value = 119123.75
$ws.Range("N135").Value = -129263.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1416.2941
$ws.Range("I107").Value = 1729
$ws.Range("K107").Value = 5187
$ws.Range("M107").Value = -3267
$ws.Range("H122").Value = 993.08826
$ws.Range("I122").Value = 612.7778
$ws.Range("J122").Value = 1130
$ws.Range("K122").Value = 5515.000199999999
$ws.Range("L122").Value = 10170
$ws.Range("M122").Value = -3065.000199999999
$ws.Range("N122").Value = -15070
$ws.Range("H132").Value = 1985.32
$ws.Range("I132").Value = 1876.8182
$ws.Range("K132").Value = 16891.3638
$ws.Range("M132").Value = -14361.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6336.864
$ws.Range("I80").Value = 4339.3335
$ws.Range("J80").Value = 10617.286
$ws.Range("K80").Value = 4339.3335
$ws.Range("L80").Value = 10617.286
$ws.Range("M80").Value = -3341.3335
$ws.Range("N80").Value = -12613.286
$ws.Range("H83").Value = 6336.864
$ws.Range("I83").Value = 4339.3335
$ws.Range("J83").Value = 10617.286
$ws.Range("K83").Value = 21696.6675
$ws.Range("L83").Value = 53086.43
$ws.Range("M83").Value = -16704.6675
$ws.Range("N83").Value = -63070.43
$ws.Range("H97").Value = 1068.3784
$ws.Range("I97").Value = 849.7143
$ws.Range("K97").Value = 849.7143
$ws.Range("M97").Value = -353.7143
$ws.Range("H107").Value = 924.6667
$ws.Range("I107").Value = 1836.3334
$ws.Range("J107").Value = 468.83334
$ws.Range("K107").Value = 1836.3334
$ws.Range("L107").Value = 468.83334
$ws.Range("M107").Value = 83.66660000000002
$ws.Range("N107").Value = -4308.83334
$ws.Range("H126").Value = 2462.261
$ws.Range("I126").Value = 2427.0625
$ws.Range("K126").Value = 7281.1875
$ws.Range("M126").Value = -4811.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5561.4375
$ws.Range("I46").Value = 850.5
$ws.Range("J46").Value = 7131.75
$ws.Range("K46").Value = 850.5
$ws.Range("L46").Value = 7131.75
$ws.Range("M46").Value = -662.5
$ws.Range("N46").Value = -7507.75
$ws.Range("H132").Value = 4122.3335
$ws.Range("I132").Value = 3508.2888
$ws.Range("J132").Value = 13333
$ws.Range("K132").Value = 10524.8664
$ws.Range("L132").Value = 39999
$ws.Range("M132").Value = -7994.866399999999
$ws.Range("N132").Value = -45059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 386840.53
$ws.Range("I126").Value = 2366.6667
$ws.Range("K126").Value = 7100.000100000001
$ws.Range("M126").Value = -4630.000100000001
